# Updates the cryptos list (Price / Volume(1h) columns, plus one
# Toncoin/RenderToken row swap) to match the latest scrape.
#
# Note: several Price values look like plain numbers (e.g. "0.593",
# "9.08"). Assigning those to Range.Value directly would make Excel
# auto-convert them to numeric cells and lose the exact text formatting
# (e.g. "2.70" -> 2.7). To keep them as text - exactly like the original
# cells - we prefix with a leading apostrophe to force text entry, then
# reset Style back to "Normal" so no extra per-cell formatting lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.023.51"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "3.357.54"
$ws.Range("E3").Value = "  +3.60%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'523.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "'173.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D7").Value = "'0.593"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "3.333.32"
$ws.Range("E8").Value = "  +3.35%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.607"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "'53.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.72%  "
$ws.Range("D12").Value = "'0.134"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "'9.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.919.44"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").Value = "3.377.51"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "'17.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "63.963.93"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "'11.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").Value = "'0.961"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "'374.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'4.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.23%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "'11.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").Value = "'81.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").Value = "'3.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'6.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("D28").Value = "'2.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").Value = "'11.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "'8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'28.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").Value = "'629.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'6.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.33%  "
$ws.Range("D34").Value = "'11.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").Value = "'57.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'36.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "'0.379"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.18%  "
$ws.Range("D40").Value = "0.0₃0732"
$ws.Range("E40").Value = "  +11.01%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "'2.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.65%  "
$ws.Range("D43").Value = "2.983.79"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").Value = "'0.125"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").Value = "'3.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.85%  "
$ws.Range("E46").Value = "  +3.12%  "
$ws.Range("D47").Value = "'0.0395"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "'2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'135.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.86%  "
